$d = $word.ActiveDocument
Write-Output "noop"
